$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (under the frozen header) was blank in those columns; fill every
# cell with 1, matching the author's edit.
$ws.Range("C3:J3").Value = 1
$ws.Range("M3").Value = 1

# Leave the cursor on M3, same as the recorded selection in the saved file.
$ws.Range("M3").Select()
